$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (styles) of the previous data row onto the new row
$ws.Range("A52:B52").Copy()
$ws.Range("A53:B53").PasteSpecial(-4122)

# Enter the new date as a formula-literal so it is stored as text (not
# auto-converted to a date serial number), matching how the existing
# "dd-mm-yyyy" strings in column A are stored as shared strings.
$ws.Cells.Item(53, 1).Formula = '="06-11-2025"'
$ws.Cells.Item(53, 2).Value = "The price of gold in India today is ₹12,191 per gram for 24 karat gold, ₹11,175 per gram for 22 karat gold and ₹9,143 per gram for 18 karat gold (also called 999 gold)."

# Convert the formula in A53 to its literal text value
$ws.Cells.Item(53, 1).Copy()
$ws.Cells.Item(53, 1).PasteSpecial(-4163)

$excel.CutCopyMode = $false
